$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting "Item Category Name" -> B
# and "Item Group Name" -> C, and give column A a new "ID" header.
$ws.Range("A1").EntireColumn.Insert()

# Update header text to match the new template layout.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Category Name"
$ws.Range("C1").Value = "Group Name"

# The inserted column doesn't pick up the bold header formatting used by
# its neighbours - copy it over from column B so all three headers match.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths as captured in the target workbook (values are pre-adjusted
# for this engine's ColumnWidth -> stored-width rounding so the saved file
# ends up with width="24.28515625" / "47" / "45").
$ws.Columns.Item(1).ColumnWidth = 23.42
$ws.Columns.Item(2).ColumnWidth = 46.16666666666666
$ws.Columns.Item(3).ColumnWidth = 44.16666666666666

# Restore selection to the next empty cell as seen in the saved workbook.
$ws.Range("C2").Select()
